# Loan RBI, Variable Instalments
# - Insert a new (blank) column at "N" on the "Repayment schedule" sheet,
#   shifting the old N/O/P ("Late"/heading/"Outstanding") data right to O/P/Q.
# - Give the new column roughly the same width as its left neighbour (M).
# - Make "Repayment schedule" the active sheet/tab, with cell R10 selected.

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Switch to the Repayment schedule sheet (this becomes the saved "active tab",
# and the "Input" sheet - previously active - stops being tabSelected).
$wsSchedule.Activate()

# Insert a blank column before the existing "N" column (Late/heading/Outstanding
# shift one column to the right, to O/P/Q).
$wsSchedule.Columns("N").Insert()

# The inherited column keeps roughly column M's width.
$wsSchedule.Columns("N").ColumnWidth = 9.83

# Leave the selection where the author left it after the edit.
$wsSchedule.Range("R10").Select()
